$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "creation.option1" row (row 10): the "Better quality" option no longer
# implies a price increase, so the translations are shortened. Setting the
# cells in this specific order (fr, en, de, it, es) reproduces the shared
# string table ordering used by the author's edit.
$ws.Range("C10").Value = "Meilleure qualité"
$ws.Range("B10").Value = "Better quality"
$ws.Range("D10").Value = "Bessere Bildqualität"
$ws.Range("E10").Value = "Migliore qualità dell'immagine"
$ws.Range("F10").Value = "Mejor calidad de imagen"

# Update the view/selection state left behind in the sheet (scrolled to and
# selecting F10 rather than D20).
$ws.Range("F10").Select()
